# Apply the cryptos.xlsx crypto-price-table refresh (GitHub Actions bot update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric ("37.247.01", "1.00", "0.0759", ...).
# Excel's normal Value-assignment auto-parses those into numbers/dates, which would
# silently corrupt the text (e.g. drop trailing zeros or mis-split on the extra dots).
# Force the column to Text first so every assignment below lands as a literal string,
# then restore the default "Normal" style so no stray formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.247.01"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.068.26"
$ws.Range("E3").Value = "  +3.94%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "235.96"
$ws.Range("E5").Value = "  -3.46%  "
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "57.24"
$ws.Range("E8").Value = "  +4.71%  "
$ws.Range("D9").Value = "0.379"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("D10").Value = "58.18"
$ws.Range("E10").Value = "  -2.19%  "
$ws.Range("D11").Value = "0.0759"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "0.101"
$ws.Range("E12").Value = "  +3.02%  "
$ws.Range("D13").Value = "2.376.37"
$ws.Range("E13").Value = "  +4.29%  "
$ws.Range("D14").Value = "14.50"
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("D15").Value = "21.23"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "0.774"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").Value = "5.24"
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("D18").Value = "2.011.81"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "37.461.85"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").Value = "5.93"
$ws.Range("E20").Value = "  +19.15%  "
$ws.Range("D21").Value = "68.30"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "0.0₃0810"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "223.58"
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").Value = "163.18"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").Value = "8.86"
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("D29").Value = "0.131"
$ws.Range("E29").Value = "  +6.17%  "
$ws.Range("D30").Value = "19.21"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("D33").Value = "4.44"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").Value = "0.0621"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("E35").Value = "  +7.39%  "
$ws.Range("D36").Value = "4.37"
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "5.94"
$ws.Range("E38").Value = "  +13.54%  "
$ws.Range("D39").Value = "3.33"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("D42").Value = "4.43"
$ws.Range("E42").Value = "  +24.11%  "
$ws.Range("D43").Value = "0.0954"
$ws.Range("E43").Value = "  +7.49%  "
$ws.Range("D44").Value = "1.469.93"
$ws.Range("E44").Value = "  +3.30%  "
$ws.Range("D45").Value = "94.52"
$ws.Range("E45").Value = "  +7.19%  "
$ws.Range("D46").Value = "0.0208"
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "1.14"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "16.13"
$ws.Range("E48").Value = "  +5.05%  "
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("E50").Value = "  +7.73%  "
$ws.Range("D51").Value = "2.93"
$ws.Range("E51").Value = "  +2.28%  "

# Restore default styling on column D now that the text values are safely in place.
$ws.Range("D2:D51").Style = "Normal"
